$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values that changed for this data refresh
$ws.Range("B140").Value = 11464
$ws.Range("B141").Value = 19037
$ws.Range("B144").Value = 29486

# Append new row 153 with the new monthly data point (01-08-2021)
# Force text number format first so Excel doesn't auto-convert the
# "dd-mm-yyyy"-looking string into a date serial value, then restore
# the default (Normal) style so no extra formatting is applied to the cell.
$ws.Range("A153").NumberFormat = "@"
$ws.Range("A153").Value = "01-08-2021"
$ws.Range("A153").Style = "Normal"
$ws.Range("B153").Value = 37564
